$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.322.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.871.17'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.03'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2887'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06630'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.72'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08045'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.49'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.873.22'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.150'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '272.22'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.314.37'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.17'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007753'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +6.21%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.118.03'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.320'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.223'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.406'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.44'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.98'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.374'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09924'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.375'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.468'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.085'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04709'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7023'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.705'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01889'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.651'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.334'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.81'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.962'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8454'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4177'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9999'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.56'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.274'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.099'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '928.86'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.53'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05684'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.61%  '
